# Rebuild index after corruption
# Append 10 new weekly rows (26-35) to the eth_purchases data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(46006, 102259, 318536785, 2964.18310546875,   10.47025078220676),
    @(46013, 98852,  294578960, 3006.07373046875,   10.29435269453544),
    @(46020, 44463,  131610480, 2934.538330078125,   6.904078079755683),
    @(46027, 32977,  100579850, 3226.13037109375,    5.945821750849223),
    @(46034, 24266,  76923220,  3092.3251953125,     5.100413773821187),
    @(46041, 35268,  116540860, 3186.62109375,       6.148890077690821),
    @(46048, 40302,  120099960, 2926.45703125,       6.573090512286873),
    @(46055, 41788,  141404692, 2344.356689453125,   6.69317386039803),
    @(46062, 40613,  88536340,  2103.567626953125,   6.598403184030338),
    @(46070, 45759,  91289205,  1992.1943359375,     7.003974814820833)
)

$r = 26
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item(25, 1).NumberFormat
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
